# Update "想去人数" (F) / "最低票价" (G) figures on the "展览" and
# "全部类型" worksheets to match the refreshed source data.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("G2").Value  = 66
$ws1.Range("F3").Value  = 1886
$ws1.Range("G3").Value  = 58.8
$ws1.Range("G4").Value  = 55
$ws1.Range("F6").Value  = 3186
$ws1.Range("F7").Value  = 585
$ws1.Range("F8").Value  = 604
$ws1.Range("G8").Value  = 55
$ws1.Range("F9").Value  = 291
$ws1.Range("G9").Value  = 55
$ws1.Range("F10").Value = 657
$ws1.Range("G10").Value = 60
$ws1.Range("F12").Value = 559
$ws1.Range("F18").Value = 1653
$ws1.Range("F21").Value = 623
$ws1.Range("F27").Value = 117
$ws1.Range("F32").Value = 4043
$ws1.Range("F34").Value = 784
$ws1.Range("F36").Value = 1612
$ws1.Range("F38").Value = 1911

# Sheet "全部类型" (All types) - mirrors the same events, offset by one
# extra row (演出/本地生活 entries interleaved) starting at row 13.
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("G2").Value  = 66
$ws4.Range("F3").Value  = 1886
$ws4.Range("G3").Value  = 58.8
$ws4.Range("G4").Value  = 55
$ws4.Range("F6").Value  = 3186
$ws4.Range("F7").Value  = 585
$ws4.Range("F8").Value  = 604
$ws4.Range("G8").Value  = 55
$ws4.Range("F9").Value  = 291
$ws4.Range("G9").Value  = 55
$ws4.Range("F10").Value = 657
$ws4.Range("G10").Value = 60
$ws4.Range("F12").Value = 559
$ws4.Range("F19").Value = 1653
$ws4.Range("F22").Value = 623
$ws4.Range("F28").Value = 117
$ws4.Range("F33").Value = 4043
$ws4.Range("F37").Value = 784
$ws4.Range("F39").Value = 1612
$ws4.Range("F41").Value = 1911
